# v1.5 modify (SRS ID) Column according to naming convention
#
# 1) Rename every "SRS ID" (column C) entry on the LH_SRS sheet from the old
#    "SRS-xxx" naming convention to the new "LH-SRS-xxx" naming convention
#    (rows 2-39).
# 2) Record the change on the LH_SRS_VERSION_HISTORY sheet by appending a
#    new v1.5 row.
# 3) Leave the workbook with the LH_SRS sheet active/selected (matching the
#    saved view state of the authored edit).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # LH_SRS
$ws2 = $wb.Worksheets.Item(2)   # LH_SRS_VERSION_HISTORY

# --- 1) Update the "SRS ID" column (C) naming convention ------------------
for ($r = 2; $r -le 39; $r++) {
    $cell = $ws1.Cells.Item($r, 3)
    $oldValue = $cell.Value2
    if ($oldValue) {
        $cell.Value = "LH-" + $oldValue
    }
}

# --- 2) Append the new version-history row ---------------------------------
# Copy the formatting of the previous entry (row 5) down into the new row 6
# so the new row matches the look of the existing history rows.
$ws2.Range("A5:D5").Copy()
$ws2.Range("A6:D6").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Rows.Item(6).RowHeight = 37.5

$ws2.Range("A6").Value = "v1.5"
$ws2.Range("B6").Value = "Hala Eldaly"
$ws2.Range("C6").Value = "modify (SRS ID) Column according to naming convention"
$ws2.Range("D6").Value = 45773

# --- 3) Restore view/selection state ---------------------------------------
# Select C6 on the version-history sheet first (becomes the "last touched"
# selection for that sheet), then activate LH_SRS and select C1 so it is
# the active sheet/tab when the workbook is saved.
$ws2.Range("C6").Select()

$ws1.Activate()
$ws1.Range("C1").Select()
